# Generate Report for Handback
# Replace the two source-file identifiers (and their derived artifact names/timestamps)
# that are baked into the "Overview", "zh-cn" and "de-de" worksheets with the new
# values produced by the latest handback run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# New values
# ---------------------------------------------------------------------------
$newMd1        = "ef701944-3b97-4934-aeeb-135e31bd005b.md"
$newMd1Path    = "e2e\ef701944-3b97-4934-aeeb-135e31bd005b.md"
$newMd2        = "ffff2c61959e-f720-4c00-b485-8d82c660b6c4.md"
$newMd2Path    = "e2e\ffff2c61959e-f720-4c00-b485-8d82c660b6c4.md"

$newXliffZh    = "ef701944-3b97-4934-aeeb-135e31bd005b.2188f9abf17651cb3137b5a4bf472707449b1a86.zh-cn.xlf"
$newXliffDe    = "ef701944-3b97-4934-aeeb-135e31bd005b.2188f9abf17651cb3137b5a4bf472707449b1a86.de-de.xlf"

$newGenDate       = "2016-08-17 19:01:45"
$newZhHandoffDate = "2016-08-17 19:01:39"
$newZhHandbackDate= "2016-08-17 19:01:56"
$newDeHandbackDate= "2016-08-17 19:02:11"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview.Range("A2").Value2 = $newMd1
$wsOverview.Range("B2").Value2 = $newMd1Path
$wsOverview.Range("G2").Value2 = $newGenDate

$wsOverview.Range("A3").Value2 = $newMd2
$wsOverview.Range("B3").Value2 = $newMd2Path
$wsOverview.Range("G3").Value2 = $newGenDate

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = $newMd1Path
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = $newMd2Path
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn.Range("A2").Value2 = $newMd1
$wsZhCn.Range("G2").Value2 = $newXliffZh
$wsZhCn.Range("H2").Value2 = $newZhHandoffDate
$wsZhCn.Range("I2").Value2 = $newMd1
$wsZhCn.Range("J2").Value2 = $newXliffZh
$wsZhCn.Range("K2").Value2 = $newZhHandbackDate

$wsZhCn.Range("A3").Value2 = $newMd2
$wsZhCn.Range("G3").Value2 = $newXliffZh
$wsZhCn.Range("H3").Value2 = $newZhHandoffDate
$wsZhCn.Range("I3").Value2 = $newMd2
$wsZhCn.Range("J3").Value2 = $newXliffZh
$wsZhCn.Range("K3").Value2 = $newZhHandbackDate

foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newMd1
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = $newMd1
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = $newMd2
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = $newMd2
    }
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe.Range("A2").Value2 = $newMd1
$wsDeDe.Range("G2").Value2 = $newXliffDe
$wsDeDe.Range("H2").Value2 = $newGenDate
$wsDeDe.Range("I2").Value2 = $newMd1
$wsDeDe.Range("J2").Value2 = $newXliffDe
$wsDeDe.Range("K2").Value2 = $newDeHandbackDate

$wsDeDe.Range("A3").Value2 = $newMd2
$wsDeDe.Range("G3").Value2 = $newXliffDe
$wsDeDe.Range("H3").Value2 = $newGenDate
$wsDeDe.Range("I3").Value2 = $newMd2
$wsDeDe.Range("J3").Value2 = $newXliffDe
$wsDeDe.Range("K3").Value2 = $newDeHandbackDate

foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newMd1
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = $newMd1
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = $newMd2
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = $newMd2
    }
}
